$wb = $excel.ActiveWorkbook

# --- Sheet R1: update elapsed-duration values and append a new outage row ---
$ws1 = $wb.Worksheets.Item("R1")
$ws1.Range("G2").Value = "3930:04:24"
$ws1.Range("G3").Value = "69:37:02"
$ws1.Range("G4").Value = "92:37:02"

$ws1.Range("A5").Value = ""
$ws1.Range("B5").Value = "R4"
$ws1.Range("C5").Value = ""
$ws1.Range("D5").Value = "JED0125"
$ws1.Range("E5").Value = ""
$ws1.Range("F5").Value = ""
$ws1.Range("G5").Value = ""
$ws1.Range("H5").Value = ""
$ws1.Range("I5").Value = "Generator-SG"
$ws1.Range("J5").Value = "Good+In progress"
$ws1.Range("K5").Value = ""
$ws1.Range("L5").Value = "Latis"

# --- Sheet R2: update elapsed-duration values ---
$ws2 = $wb.Worksheets.Item("R2")
$ws2.Range("G2").Value = "12111:28:05"
$ws2.Range("G3").Value = "3241:11:34"
$ws2.Range("G4").Value = "479:23:08"

# --- Sheet R4: update elapsed-duration values ---
$ws4 = $wb.Worksheets.Item("R4")
$ws4.Range("G2").Value = "2957:17:54"
$ws4.Range("G3").Value = "184:30:09"
$ws4.Range("G4").Value = "72:42:34"
$ws4.Range("G5").Value = "70:20:07"

# --- Sheet R5: update elapsed-duration value ---
$ws5 = $wb.Worksheets.Item("R5")
$ws5.Range("G2").Value = "431:16:53"

# --- Sheet R6: update elapsed-duration value ---
$ws6 = $wb.Worksheets.Item("R6")
$ws6.Range("G2").Value = "71:49:11"
